$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 261, shifting existing rows 261:274 down to 262:275
$ws.Rows.Item(261).Insert()

# Populate the newly inserted row 261 with the new weekly record
$ws.Cells.Item(261, 1).Value = 8
$ws.Cells.Item(261, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(261, 3).Value = "Coquimbo"
$ws.Cells.Item(261, 4).Value = 45267
$ws.Cells.Item(261, 5).Value = 4
$ws.Cells.Item(261, 6).Value = 100112040
$ws.Cells.Item(261, 7).Value = "Cilantro"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 2000
$ws.Cells.Item(261, 11).Value = 2300
$ws.Cells.Item(261, 12).Value = 2500
$ws.Cells.Item(261, 13).Value = 2400
$ws.Cells.Item(261, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(261, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(261, 16).Value = 1600
$ws.Cells.Item(261, 17).Value = 1.5
$ws.Cells.Item(261, 18).Value = "Hortaliza"
